$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.185.50"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.585.71"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.47"
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("E6").Value = "  +0.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.245"
$ws.Range("E8").Value = "  +0.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0606"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.24"
$ws.Range("E10").Value = "  -1.60%  "

$ws.Range("E11").Value = "  +0.66%  "

$ws.Range("D12").Value = "1.809.12"
$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").Value = "1.599.46"
$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("E14").Value = "  -1.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  +0.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.04"
$ws.Range("E16").Value = "  -0.55%  "

$ws.Range("D17").Value = "26.177.17"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").Value = "0.0₃0724"
$ws.Range("E18").Value = "  -0.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.36"
$ws.Range("E19").Value = "  +1.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "212.35"
$ws.Range("E20").Value = "  +1.74%  "

$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.97"
$ws.Range("E24").Value = "  +1.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.51"
$ws.Range("E25").Value = "  -0.27%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("E28").Value = "  -0.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.14"
$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("E30").Value = "  -1.90%  "

$ws.Range("E31").Value = "  +1.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("E32").Value = "  -0.75%  "

$ws.Range("D33").Value = "1.335.86"
$ws.Range("E33").Value = "  +4.31%  "

$ws.Range("E34").Value = "  -1.99%  "

$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.45"
$ws.Range("E36").Value = "  -1.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.579"
$ws.Range("E37").Value = "  -5.07%  "

$ws.Range("E38").Value = "  -0.56%  "

$ws.Range("E39").Value = "  +1.03%  "

$ws.Range("E40").Value = "  +3.67%  "

$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.970"
$ws.Range("E42").Value = "  -13.06%  "

$ws.Range("E43").Value = "  +0.83%  "

$ws.Range("E44").Value = "  +0.54%  "

$ws.Range("D45").Value = "1.721.27"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.93"
$ws.Range("E46").Value = "  -2.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.82"
$ws.Range("E47").Value = "  -3.22%  "

$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -0.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.47"
$ws.Range("E49").Value = "  -1.62%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0983"
$ws.Range("E50").Value = "  -2.12%  "

$ws.Range("E51").Value = "  -0.82%  "
